{"js": "// Add Bophelo's student number to the SRS.\n// The document currently has a placeholder \"(?)\" right after\n// \"Bophelo Malepu\" (team member list). Replace that placeholder with\n// her actual student number \"(20232645)\", keeping the surrounding\n// bold formatting intact.\n\nconst body = context.document.body;\n\n// Find the placeholder text. It only occurs once in the whole document,\n// right after \"Bophelo Malepu\" in the team-members section.\nconst results = body.search(\"(?)\", {\n  matchCase: true,\n  matchWholeWord: false,\n  matchWildcards: false,\n});\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the \"(?)\" placeholder to replace.');\n}\n\n// Replace the placeholder in place so the run's existing character\n// formatting (bold) carries over to the new text.\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"(20232645)\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Add Bophelo's student number to the SRS.\n# The document currently has a placeholder \"(?)\" right after\n# \"Bophelo Malepu\" (team member list). Replace that placeholder with\n# her actual student number \"(20232645)\", keeping the surrounding\n# bold formatting intact.\n\n$d = $word.ActiveDocument\n\n# Locate the lone \"(?)\" placeholder in the document and replace just\n# that range's text, so the run's existing formatting (bold) is kept.\n$range = $d.Content\n$range.Find.MatchWildcards = $false\n$range.Find.MatchCase = $true\n$found = $range.Find.Execute(\"(?)\")\n\nif ($found) {\n    $range.Text = \"(20232645)\"\n}\n"}
